$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, D hold plain-text values in the source data (some numeric-
# looking, e.g. "250.96" or "37.319.87"). Force text format before writing so
# Excel does not silently reinterpret them as numbers, then clear the format
# again so no stray cell style is left behind (matches original unstyled cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.319.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.091.61"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.96"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.73%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.23"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +22.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.74"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0742"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("E12").Value = "  +8.58%  "
$ws.Range("E13").Value = "  +6.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.398.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.52%  "
$ws.Range("E15").Value = "  +4.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.107.01"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.25%  "
$ws.Range("E17").Value = "  +6.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.225.29"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +14.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0846"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.68"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.88%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.25"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.58%  "
$ws.Range("E29").Value = "  +3.60%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.61%  "
$ws.Range("E32").Value = "  +28.30%  "
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0616"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0919"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.17%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  +7.47%  "
$ws.Range("E40").Value = "  +2.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.27"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +15.66%  "
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("E43").Value = "  +6.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.65"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0935"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +15.27%  "
$ws.Range("E46").Value = "  +114.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.322.54"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +6.16%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "RenderToken"
$ws.Range("B50").ClearFormats()
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C50").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.52%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "FraxShare"
$ws.Range("B51").ClearFormats()
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C51").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.99"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +15.28%  "
